$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet and the workbook's tab entry -------------------------
$ws.Name = "Ceny předplacených karet"

# --- Add the 10 new "Ha-loo mobile" rows (35-44) ----------------------------
# Columns: A=Spolecnost, B=Nazev produktu, C=Pocet dat, D=Cena za 1 MB (formula),
#          E=Cena, F=Mesicne, G=Nejmensi jednotka

$rows = @(
    @{ Row=35; Produkt="Data pro tablet"; Pocet=3000;  Cena=440 },
    @{ Row=36; Produkt="Data pro tablet"; Pocet=10000; Cena=640 },
    @{ Row=37; Produkt="Data pro mobil";  Pocet=50;    Cena=35  },
    @{ Row=38; Produkt="Data pro mobil";  Pocet=100;   Cena=60  },
    @{ Row=39; Produkt="Data pro mobil";  Pocet=150;   Cena=69  },
    @{ Row=40; Produkt="Data pro mobil";  Pocet=200;   Cena=110 },
    @{ Row=41; Produkt="Data pro mobil";  Pocet=300;   Cena=129 },
    @{ Row=42; Produkt="Data pro mobil";  Pocet=600;   Cena=199 },
    @{ Row=43; Produkt="Data pro mobil";  Pocet=1200;  Cena=289 },
    @{ Row=44; Produkt="Data pro mobil";  Pocet=1500;  Cena=349 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Range("A$i").Value = "Ha-loo mobile"
    $ws.Range("B$i").Value = $r.Produkt
    $ws.Range("C$i").Value = $r.Pocet
    $ws.Range("D$i").Formula = "=E$i/C$i"
    $ws.Range("E$i").Value = $r.Cena
    $ws.Range("F$i").Value = "ANO"
    $ws.Range("G$i").Value = "1 kB"
}

# Rows 36-44 of column D already carried a leftover placeholder style from
# blank filler rows further down the sheet; re-apply the same number format
# used by the rest of the "Cena za 1 MB" column (copied from D34) so the new
# cells match their neighbours instead of keeping that stale formatting.
$ws.Range("D34").Copy() | Out-Null
$ws.Range("D35:D44").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Hyperlinks on the "Nazev produktu" column for the new rows ------------
# Mirrors the source structure: three single-cell links, then two links that
# cover overlapping multi-cell ranges (B38:B43 and B43:B44).
$ws.Hyperlinks.Add($ws.Range("B35"), "https://www.ha-loo.cz/tarify/data-pro-tablet") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B36"), "https://www.ha-loo.cz/tarify/data-pro-tablet") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B37"), "https://www.ha-loo.cz/tarify/data-pro-mobil") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B38:B43"), "https://www.ha-loo.cz/tarify/data-pro-mobil", "", "", "Data pro mobil") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B43:B44"), "https://www.ha-loo.cz/tarify/data-pro-mobil", "", "", "Data pro mobil") | Out-Null

# Re-apply the same cell formatting (font/underline) used by the other
# hyperlink cells in column B so the new cells pick up the existing
# "Hypertextovy odkaz" style instead of a freshly minted duplicate.
$ws.Range("B33").Copy() | Out-Null
$ws.Range("B35:B44").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Extend the Excel table (ListObject) to cover the new rows -------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G44"))

# --- Update the view state (selection / scroll position) -------------------
$ws.Range("I10").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
